# Auto-generated edit script for Phantom_Profits (FFXIV leve-profit tracker)
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 369.46667
$ws.Range("I2").Value = 387.08334
$ws.Range("K2").Value = 387.08334
$ws.Range("M2").Value = -274.08334
$ws.Range("H18").Value = 4595.8335
$ws.Range("I18").Value = 4805.5557
$ws.Range("J18").Value = 3966.6667
$ws.Range("K18").Value = 4805.5557
$ws.Range("L18").Value = 3966.6667
$ws.Range("M18").Value = -4521.5557
$ws.Range("N18").Value = -4534.6667
$ws.Range("H21").Value = 7573.7144
$ws.Range("J21").Value = 7666.6665
$ws.Range("L21").Value = 7666.6665
$ws.Range("N21").Value = -8602.666499999999
$ws.Range("H23").Value = 7573.7144
$ws.Range("J23").Value = 7666.6665
$ws.Range("L23").Value = 7666.6665
$ws.Range("N23").Value = -8134.6665
$ws.Range("H40").Value = 3200
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("H53").Value = 917.6
$ws.Range("J53").Value = 461.7143
$ws.Range("L53").Value = 461.7143
$ws.Range("N53").Value = -1735.7143
$ws.Range("H55").Value = 184.83333
$ws.Range("I55").Value = 191.8
$ws.Range("J55").Value = 150
$ws.Range("K55").Value = 191.8
$ws.Range("L55").Value = 150
$ws.Range("M55").Value = 22.19999999999999
$ws.Range("N55").Value = -578
$ws.Range("H74").Value = 8248.6
$ws.Range("I74").Value = 8248.6
$ws.Range("K74").Value = 8248.6
$ws.Range("M74").Value = -7312.6
$ws.Range("H77").Value = 8248.6
$ws.Range("I77").Value = 8248.6
$ws.Range("K77").Value = 41243
$ws.Range("M77").Value = -36563
$ws.Range("H137").Value = 37038650
$ws.Range("I137").Value = 83334456
$ws.Range("K137").Value = 250003368
$ws.Range("M137").Value = -250000818
$ws.Range("N40").ClearContents()

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4507.3667
$ws.Range("I32").Value = 3758.0715
$ws.Range("K32").Value = 3758.0715
$ws.Range("M32").Value = -3471.0715
$ws.Range("H61").Value = 2492.5
$ws.Range("I61").Value = 2492.5
$ws.Range("K61").Value = 2492.5
$ws.Range("M61").Value = -2280.5
$ws.Range("H62").Value = 40226
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("H65").Value = 40226
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("H136").Value = 2492.5
$ws.Range("I136").Value = 2492.5
$ws.Range("K136").Value = 7477.5
$ws.Range("M136").Value = -4927.5
$ws.Range("N62").ClearContents()
$ws.Range("N65").ClearContents()

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1428.875
$ws.Range("I94").Value = 1139.2727
$ws.Range("J94").Value = 2066
$ws.Range("K94").Value = 1139.2727
$ws.Range("L94").Value = 2066
$ws.Range("M94").Value = -688.2727
$ws.Range("N94").Value = -2968
$ws.Range("H107").Value = 1247.5
$ws.Range("I107").Value = 1247.5
$ws.Range("K107").Value = 1247.5
$ws.Range("M107").Value = 672.5

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1754.1428
$ws.Range("I122").Value = 1687.5454
$ws.Range("K122").Value = 5062.6362
$ws.Range("M122").Value = -2612.6362
$ws.Range("H134").Value = 1336.0714
$ws.Range("I134").Value = 1390.0834
$ws.Range("K134").Value = 4170.2502
$ws.Range("M134").Value = -1635.2502

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 412.5
$ws.Range("H121").Value = 2784.4614
$ws.Range("I121").Value = 456.85715
$ws.Range("J121").Value = 5500
$ws.Range("K121").Value = 1370.57145
$ws.Range("L121").Value = 16500
$ws.Range("M121").Value = -60.57144999999991
$ws.Range("N121").Value = -19120

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 21540.75
$ws.Range("J43").Value = 62625
$ws.Range("L43").Value = 62625
$ws.Range("N43").Value = -62927
$ws.Range("H57").Value = 24599.8
$ws.Range("I57").Value = 3750
$ws.Range("J57").Value = 38499.668
$ws.Range("K57").Value = 3750
$ws.Range("L57").Value = 38499.668
$ws.Range("M57").Value = -2930
$ws.Range("N57").Value = -40139.668
$ws.Range("H102").Value = 2819.6924
$ws.Range("I102").Value = 2221.3333
$ws.Range("K102").Value = 2221.3333
$ws.Range("M102").Value = -599.3332999999998
$ws.Range("H136").Value = 39536.04
$ws.Range("I136").Value = 200000
$ws.Range("J136").Value = 32850.043
$ws.Range("K136").Value = 600000
$ws.Range("L136").Value = 98550.12899999999
$ws.Range("M136").Value = -597450
$ws.Range("N136").Value = -103650.129

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 13999.167
$ws.Range("J14").Value = 13999.167
$ws.Range("L14").Value = 13999.167
$ws.Range("N14").Value = -14343.167
$ws.Range("H22").Value = 21739692
$ws.Range("I22").Value = 434.84616
$ws.Range("J22").Value = 50000724
$ws.Range("K22").Value = 434.84616
$ws.Range("L22").Value = 50000724
$ws.Range("M22").Value = -139.84616
$ws.Range("N22").Value = -50001314
$ws.Range("H27").Value = 21739692
$ws.Range("I27").Value = 434.84616
$ws.Range("J27").Value = 50000724
$ws.Range("K27").Value = 434.84616
$ws.Range("L27").Value = 50000724
$ws.Range("M27").Value = -327.84616
$ws.Range("N27").Value = -50000938
$ws.Range("H39").Value = 9400
$ws.Range("I39").Value = 2000
$ws.Range("J39").Value = 11250
$ws.Range("K39").Value = 2000
$ws.Range("L39").Value = 11250
$ws.Range("M39").Value = -1540
$ws.Range("N39").Value = -12170
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("H46").Value = 3686.5293
$ws.Range("I46").Value = 1617
$ws.Range("J46").Value = 5135.2
$ws.Range("K46").Value = 1617
$ws.Range("L46").Value = 5135.2
$ws.Range("M46").Value = -1429
$ws.Range("N46").Value = -5511.2
$ws.Range("H48").Value = 37498.75
$ws.Range("I48").Value = 37498.75
$ws.Range("K48").Value = 37498.75
$ws.Range("M48").Value = -36837.75
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("H55").Value = 858.375
$ws.Range("I55").Value = 431.15384
$ws.Range("J55").Value = 1363.2727
$ws.Range("K55").Value = 431.15384
$ws.Range("L55").Value = 1363.2727
$ws.Range("M55").Value = -258.15384
$ws.Range("N55").Value = -1709.2727
$ws.Range("H74").Value = 25000
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("H76").Value = 16650
$ws.Range("J76").Value = 16650
$ws.Range("L76").Value = 16650
$ws.Range("N76").Value = -17326
$ws.Range("H77").Value = 25000
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("H79").Value = 16650
$ws.Range("J79").Value = 16650
$ws.Range("L79").Value = 16650
$ws.Range("N79").Value = -18990
$ws.Range("H93").Value = 718.44446
$ws.Range("I93").Value = 672
$ws.Range("J93").Value = 1090
$ws.Range("K93").Value = 672
$ws.Range("L93").Value = 1090
$ws.Range("M93").Value = 576
$ws.Range("N93").Value = -3586
$ws.Range("N42").ClearContents()
$ws.Range("N49").ClearContents()
$ws.Range("N74").ClearContents()
$ws.Range("N77").ClearContents()

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("H39").Value = 25000
$ws.Range("J39").Value = 25000
$ws.Range("L39").Value = 25000
$ws.Range("N39").Value = -25826
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("H44").Value = 32500
$ws.Range("J44").Value = 35000
$ws.Range("L44").Value = 35000
$ws.Range("N44").Value = -36108
$ws.Range("H54").Value = 47680.875
$ws.Range("J54").Value = 47680.875
$ws.Range("L54").Value = 47680.875
$ws.Range("N54").Value = -48720.875
$ws.Range("H107").Value = 942.8461
$ws.Range("I107").Value = 1061.091
$ws.Range("J107").Value = 856.13336
$ws.Range("K107").Value = 3183.273
$ws.Range("L107").Value = 2568.40008
$ws.Range("M107").Value = -1263.273
$ws.Range("N107").Value = -6408.40008
$ws.Range("H124").Value = 26621.125
$ws.Range("J124").Value = 26621.125
$ws.Range("L124").Value = 26621.125
$ws.Range("N37").ClearContents()
$ws.Range("N40").ClearContents()
